# Kaldheim Commander (KHC) sheet gains two new card rows:
#   "Inspired Sphinx" is inserted before "Lathril, Blade of the Elves"
#   "Wolverine Riders" is appended after "Ranar the Ever-Watchful"
# Write bottom-up so each source value is set before it would be overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Wolverine Riders"
$ws.Range("A4").Value = "Ranar the Ever-Watchful"
$ws.Range("A3").Value = "Lathril, Blade of the Elves"
$ws.Range("A2").Value = "Inspired Sphinx"
